$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (F2:F11) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 8987
$ws1.Range("F3").Value = 98
$ws1.Range("F4").Value = 239
$ws1.Range("F5").Value = 106
$ws1.Range("F6").Value = 1495
$ws1.Range("F7").Value = 1414
$ws1.Range("F8").Value = 250
$ws1.Range("F9").Value = 47
$ws1.Range("F10").Value = 329
$ws1.Range("F11").Value = 91

# --- Sheet "全部类型" (F2:F8, F10:F12) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 8987
$ws4.Range("F3").Value = 98
$ws4.Range("F4").Value = 239
$ws4.Range("F5").Value = 106
$ws4.Range("F6").Value = 1495
$ws4.Range("F7").Value = 1414
$ws4.Range("F8").Value = 250
$ws4.Range("F10").Value = 47
$ws4.Range("F11").Value = 329
$ws4.Range("F12").Value = 91
